$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.378.11'
$ws.Range('E2').Value = '  -4.26%  '
$ws.Range('D3').Value = '3.869.13'
$ws.Range('E3').Value = '  -4.34%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.73'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.99%  '
$ws.Range('D7').Value = '4.026.18'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.669'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.998'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.717'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.161'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '52.49'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +11.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000305'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.36%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.05%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '4.494.83'
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('D16').Value = '3.889.06'
$ws.Range('E16').Value = '  -3.54%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.131'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.88%  '
$ws.Range('D21').Value = '69.446.58'
$ws.Range('E21').Value = '  -3.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '416.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '93.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +14.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.41'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '674.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.89'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '47.15'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +11.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.124'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '66.87'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.417'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.34'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.62%  '
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.144'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.55%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0777'
$ws.Range('E43').Value = '  -9.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0465'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.144'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.22%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '27.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.49%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.69%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000266'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.10%  '
